# çıktıyı farklı bir sayfa üzerine kayıt işlemi eklendi
# Move the SORULAR/CEVAP table onto its own "BOT_CEVAPLAR" sheet, keep only the
# question + target-report columns on "Sayfa1", and refresh the two answer texts.

$wb = $excel.ActiveWorkbook
$sayfa1 = $wb.Worksheets.Item(1)

# 1) Duplicate Sayfa1 (with all five columns + formatting) right after itself,
#    then rename the duplicate to BOT_CEVAPLAR - this becomes the full answer log.
$sayfa1.Copy($null, $sayfa1)
$bot = $wb.Worksheets.Item(2)
$bot.Name = "BOT_CEVAPLAR"

# 2) Refresh the CEVAP (answer) + GİTTİĞİ RAPOR values on the new sheet with the latest text.
$gittigiRapor = "`n                    📁 Günlük POS İşlemleri ve Banka Bazlı Toplamlar`n          "
$bot.Range("C2").Value = $gittigiRapor
$bot.Range("C3").Value = $gittigiRapor
$bot.Range("E2").Value = "Bugün en yüksek POS girişi Ziraat Bankası'ndan, toplam 614,790.30 TL olarak gerçekleşmiştir."
$bot.Range("E3").Value = "Bugünkü banka bazlı girişler ve çıkışlar farkı en yüksek banka `"AKBANK POS HS.`" olup, fark tutarı -4,276,583.43 TL'dir."

# 3) Trim Sayfa1 back down to just the SORULAR / GİTMESİ GEREKEN RAPOR columns.
$sayfa1.Range("C1:E3").ClearContents()

# 4) Column widths, close to the authored layout.
$sayfa1.Columns.Item(1).ColumnWidth = 67.625
$sayfa1.Columns.Item(2).ColumnWidth = 57.875
$sayfa1.Columns.Item(3).ColumnWidth = 34.625

$bot.Columns.Item(2).ColumnWidth = 31.375
$bot.Columns.Item(3).ColumnWidth = 37.75

# 5) Leave Sayfa1 as the active/selected sheet, matching the saved view state.
$sayfa1.Activate()
$sayfa1.Range("A13").Select()
